$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("mySheet")

# Update the total-time counter values (row 2: hours/mins/secs) with the
# latest numbers, keeping them stored as text (matching the
# "numberStoredAsText" ignored-error already present on A1:C2).
$ws.Range("A2:C2").NumberFormat = "@"
$ws.Range("A2").Value = "221"
$ws.Range("B2").Value = "20"
$ws.Range("C2").Value = "48"
